# Proof of concept complete - It's working!
#
# 1. "Surnames" sheet gets renamed to "Prez" and filled with the list of
#    US presidents (column A) plus a `=""""&A#&""","` helper column (B).
# 2. "Prez" becomes the active sheet/tab.
# 3. The previously-active "Afflictions" sheet keeps its scrolled
#    position but its selection moves to B34 (a single cell, not the old
#    B1:B63 block) and loses tabSelected.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Afflictions: move the selection, single cell B34.
# ---------------------------------------------------------------------
$wsAfflictions = $wb.Worksheets.Item("Afflictions")
$wsAfflictions.Activate()
[void]$wsAfflictions.Range("B34").Select()

# ---------------------------------------------------------------------
# 2. Surnames -> Prez, filled with presidents.
# ---------------------------------------------------------------------
$wsPrez = $wb.Worksheets.Item("Surnames")
$wsPrez.Name = "Prez"

$presidents = @(
    "Washington", "Adams",     "Jefferson", "Madison",   "Monroe",
    "Adams",      "Jackson",   "Van Buren", "Harrison",  "Tyler",
    "Polk",       "Taylor",    "Fillmore",  "Pierce",    "Buchanan",
    "Lincoln",    "Johnson",   "Grant",     "Hayes",     "Garfield",
    "Arthur",     "Cleveland", "Harrison",  "Cleveland", "McKinley",
    "Roosevelt",  "Taft",      "Wilson",    "Harding",   "Coolidge",
    "Hoover",     "Roosevelt", "Truman",    "Eisenhower","Kennedy",
    "Johnson",    "Nixon",     "Ford",      "Carter",    "Reagan",
    "Bush",       "Clinton",   "Bush",      "Obama",     "Trump"
)

# Row 8 is "Van Buren" -- the 8th president, inaugurated between Jackson
# and Harrison. Every other name is written first (in row order) and
# "Van Buren" is filled in last, which is what reproduces the shared
# string table order of the target workbook (Van Buren ends up as the
# very last new shared string instead of the 7th).
for ($i = 0; $i -lt $presidents.Length; $i++) {
    if ($i -eq 7) { continue }
    $wsPrez.Cells.Item($i + 1, 1).Value = $presidents[$i]
}
$wsPrez.Cells.Item(8, 1).Value = $presidents[7]

$wsPrez.Range("B1:B45").Formula = '=""""&A1&""","'

# Auto-size columns A/B to fit the longest entries (matches the diff's
# <col .../> bestFit widths).
[void]$wsPrez.Range("A1:B45").EntireColumn.AutoFit()

$wsPrez.Activate()
